$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- 1. Set cell values (text / numbers / blanks) ----
$ws.Range("A1").Value = "Week Ending In:"
$ws.Range("B1").Value = 41019
$ws.Range("C1").Value = 41026
$ws.Range("D1").Value = 41033
$ws.Range("E1").Value = 41040
$ws.Range("F1").Value = 41047
$ws.Range("G1").Value = 41054
$ws.Range("H1").Value = 41061
$ws.Range("I1").Value = 41068
$ws.Range("J1").Value = 41075
$ws.Range("K1").Value = 41082
$ws.Range("L1").Value = 41089
$ws.Range("M1").Value = 41096
$ws.Range("N1").Value = 41103
$ws.Range("O1").Value = 41110
$ws.Range("P1").Value = 41117
$ws.Range("Q1").Value = 41124
$ws.Range("R1").Value = 41131
$ws.Range("S1").Value = 41138
$ws.Range("T1").Value = 41145
$ws.Range("A2").Value = "Mechanical/ Electrical"
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = "Fabricate Rails"
$ws.Range("E2").Value = "Mount All Electronics"
$ws.Range("F2").Value = "Wiring"
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("P2").ClearContents()
$ws.Range("Q2").ClearContents()
$ws.Range("R2").ClearContents()
$ws.Range("S2").ClearContents()
$ws.Range("T2").ClearContents()
$ws.Range("A3").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = "Assemble Frame"
$ws.Range("E3").Value = "Order Parts"
$ws.Range("F3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("P3").ClearContents()
$ws.Range("Q3").ClearContents()
$ws.Range("R3").ClearContents()
$ws.Range("S3").ClearContents()
$ws.Range("T3").ClearContents()
$ws.Range("A4").Value = "ROS Gateway"
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = "Get IRB 120 Model in Simulator"
$ws.Range("G4").Value = "ROS Arm Control"
$ws.Range("H4").ClearContents()
$ws.Range("I4").Value = "Reverse-Engineer ABB Protocol"
$ws.Range("J4").Value = "ROS Gateway"
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("P4").ClearContents()
$ws.Range("Q4").ClearContents()
$ws.Range("R4").ClearContents()
$ws.Range("S4").ClearContents()
$ws.Range("T4").ClearContents()
$ws.Range("A5").Value = "Mobility"
$ws.Range("B5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("I5").Value = "Driving by Remote Control (ROS Bringup)"
$ws.Range("J5").Value = "Reflexive Halt +"
$ws.Range("K5").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("P5").ClearContents()
$ws.Range("Q5").ClearContents()
$ws.Range("R5").ClearContents()
$ws.Range("S5").ClearContents()
$ws.Range("T5").ClearContents()
$ws.Range("A6").Value = "Rapid"
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("J6").Value = "Rapid Stub"
$ws.Range("K6").ClearContents()
$ws.Range("L6").ClearContents()
